$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and date range) ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Data table updates (rows 15-30) ---
$ws.Range("D15").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("G15").Value = 2
$ws.Range("N15").Value = -46.153846153846
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("I15").Copy()
$ws.Range("D16").PasteSpecial(-4122)
$ws.Range("E16").Value = 25
$ws.Range("H15").Copy()
$ws.Range("E16").PasteSpecial(-4122)
$ws.Range("F16").Value = 11
$ws.Range("G16").Value = 12
$ws.Range("H16").Value = -8.333333333333
$ws.Range("I16").Value = 49
$ws.Range("J16").Value = 63
$ws.Range("K16").Value = -22.222222222222
$ws.Range("L16").Value = 75
$ws.Range("M16").Value = -45.555555555555
$ws.Range("N16").Value = -84.444444444444
$ws.Range("C17").Value = 1
$ws.Range("E17").Value = -80
$ws.Range("F17").Value = 14
$ws.Range("G17").Value = 15
$ws.Range("H17").Value = -6.666666666666
$ws.Range("J17").Value = 55
$ws.Range("K17").Value = 40
$ws.Range("M17").Value = 140.625
$ws.Range("N17").Value = 1.315789473684
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 20
$ws.Range("F18").Value = 21
$ws.Range("G18").Value = 14
$ws.Range("H18").Value = 50
$ws.Range("I18").Value = 95
$ws.Range("J18").Value = 60
$ws.Range("K18").Value = 58.333333333333
$ws.Range("L18").Value = 72.727272727272
$ws.Range("M18").Value = 11.764705882352
$ws.Range("N18").Value = -81.904761904761
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 14
$ws.Range("E19").Value = -14.285714285714
$ws.Range("F19").Value = 41
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -4.651162790697
$ws.Range("I19").Value = 208
$ws.Range("J19").Value = 222
$ws.Range("K19").Value = -6.306306306306
$ws.Range("L19").Value = 114.432989690722
$ws.Range("M19").Value = 26.829268292682
$ws.Range("N19").Value = 1.463414634146
$ws.Range("C20").Value = 4
$ws.Range("E20").Value = 33.333333333333
$ws.Range("F20").Value = 19
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 26.666666666666
$ws.Range("I20").Value = 89
$ws.Range("J20").Value = 85
$ws.Range("K20").Value = 4.705882352941
$ws.Range("L20").Value = 217.857142857143
$ws.Range("M20").Value = 25.352112676056
$ws.Range("N20").Value = -94.825581395348
$ws.Range("C21").Value = 28
$ws.Range("D21").Value = 31
$ws.Range("E21").Value = -9.677419354838
$ws.Range("F21").Value = 106
$ws.Range("H21").Value = 4.950495049504
$ws.Range("I21").Value = 525
$ws.Range("J21").Value = 495
$ws.Range("K21").Value = 6.060606060606
$ws.Range("L21").Value = 109.163346613546
$ws.Range("M21").Value = 17.977528089887
$ws.Range("N21").Value = -81.617647058823
$ws.Range("C22").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("F22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("L22").Value = 33.333333333333
$ws.Range("C23").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C23").PasteSpecial(-4122)
$ws.Range("D23").Value = 2
$ws.Range("E23").Value = -100
$ws.Range("F23").Value = 6
$ws.Range("G23").Value = 5
$ws.Range("H23").Value = 20
$ws.Range("I23").Value = 29
$ws.Range("J23").Value = 21
$ws.Range("K23").Value = 38.095238095238
$ws.Range("L23").Value = 141.666666666667
$ws.Range("M23").Value = 163.636363636364
$ws.Range("C24").Value = 33
$ws.Range("D24").Value = 48
$ws.Range("E24").Value = -31.25
$ws.Range("F24").Value = 120
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -9.090909090909
$ws.Range("I24").Value = 529
$ws.Range("J24").Value = 438
$ws.Range("K24").Value = 20.776255707762
$ws.Range("L24").Value = 104.247104247104
$ws.Range("M24").Value = 79.931972789115
$ws.Range("C25").Value = 5
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -28.571428571428
$ws.Range("F25").Value = 30
$ws.Range("H25").Value = -14.285714285714
$ws.Range("I25").Value = 160
$ws.Range("J25").Value = 143
$ws.Range("K25").Value = 11.888111888111
$ws.Range("L25").Value = 39.130434782608
$ws.Range("M25").Value = 7.382550335570
$ws.Range("C26").Value = 1
$ws.Range("I15").Copy()
$ws.Range("C26").PasteSpecial(-4122)
$ws.Range("E26").Value = 0
$ws.Range("G26").Value = 3
$ws.Range("H26").Value = -66.666666666666
$ws.Range("I26").Value = 11
$ws.Range("J26").Value = 14
$ws.Range("K26").Value = -21.428571428571
$ws.Range("L26").Value = 83.333333333333
$ws.Range("C27").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 3
$ws.Range("I15").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("H15").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 50
$ws.Range("J27").Value = 21
$ws.Range("K27").Value = -19.047619047619
$ws.Range("L27").Value = 54.545454545454
$ws.Range("D28").Value = 2
$ws.Range("G28").Value = 6
$ws.Range("J28").Value = 6
$ws.Range("K28").Value = -83.333333333333
$ws.Range("G29").Value = 4
$ws.Range("J29").Value = 4
$ws.Range("K29").Value = -75
$ws.Range("D30").Value = "'0"
$ws.Range("C14").Copy()
$ws.Range("D30").PasteSpecial(-4122)
$ws.Range("E30").Value = "'***.*"
$ws.Range("E14").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("I30").Value = 9
$ws.Range("K30").Value = 350
